$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

$wsZh.Range("D4").Value = "2016-01-19 06:32:03"
$wsDe.Range("D4").Value = "2016-01-19 06:32:13"
